# "update icons for plot styles"
#
# On slide 3, the 13 shapes that sit next to "Group 66" in the plot-style
# icon (the original "Group 66" connector-group plus 12 "Straight
# Connector" shapes that used to be loose siblings in the spTree) get
# wrapped in one new outer group. "Group 66" keeps its own original
# child-offset coordinate space as a nested sub-group, while the new
# outer group's child space lines up 1:1 with its own position (the
# normal state for a freshly created PowerPoint group).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# PowerPoint never reuses shape IDs within a session, so the new group
# that was created here ended up as id 120 ("Group 119") even though the
# highest id already present on the slide is 112. Reproduce that by
# burning through the intermediate, already-taken ids with harmless
# add+delete cycles before doing the real grouping operation, so the new
# group lands on the same id/name as the authored edit.
for ($i = 0; $i -lt 80; $i++) {
    $tmp = $s.Shapes.AddShape(1, 0, 0, 10, 10)
    $tmp.Delete()
}

$names = @(
    "Group 66",
    "Straight Connector 80",
    "Straight Connector 81",
    "Straight Connector 83",
    "Straight Connector 85",
    "Straight Connector 86",
    "Straight Connector 87",
    "Straight Connector 103",
    "Straight Connector 104",
    "Straight Connector 105",
    "Straight Connector 109",
    "Straight Connector 110",
    "Straight Connector 111"
)
$range = $s.Shapes.Range($names)
$newGroup = $range.Group()
$newGroup.Name = "Group 119"
